$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column width adjustments (columns A-C shrink a bit, E-F adjust, and the
#    previously-merged G:H width group is split into its own G and H widths)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 38.5546875
$ws.Columns.Item(2).ColumnWidth = 30.21875
$ws.Columns.Item(3).ColumnWidth = 19.6640625
$ws.Columns.Item(5).ColumnWidth = 29.88671875
$ws.Columns.Item(6).ColumnWidth = 29.33203125
$ws.Columns.Item(7).ColumnWidth = 23.33203125
$ws.Columns.Item(8).ColumnWidth = 25.21875

# ---------------------------------------------------------------------------
# 2. Clear the (essentially invisible) extra style that had been applied to
#    column H in rows 62-85. This also makes now-empty rows (73 and 85)
#    disappear from the sheet, matching the cleaned-up sheet.
# ---------------------------------------------------------------------------
$hRows = 62,63,64,65,66,67,68,69,70,71,72,74,75,76,77,78,79,80,81,82,83,84
foreach ($r in $hRows) {
    $ws.Range("H$r").ClearFormats()
}
$ws.Range("H73").Clear()
$ws.Range("H85").Clear()

# ---------------------------------------------------------------------------
# 3. New summary columns (E, F, G) in the analysis table at the bottom of
#    the sheet (rows 86-97).
# ---------------------------------------------------------------------------
$ws.Range("E86").Value = "val before upscale"
$ws.Range("F86").Value = "val after upscale (4x lite pad=1)"
$ws.Range("G86").Value = "test before upscale"

$ws.Range("E87").Value = 0.79873551106427798
$ws.Range("F87").Value = 0.78205128205128205
$ws.Range("G87").Value = 0.66700000000000004

$ws.Range("E88").Value = 0.82033719704952501
$ws.Range("E89").Value = 0.815946610467158
$ws.Range("E90").Value = 0.87671232876712302
$ws.Range("E91").Value = 0.769582016157358
$ws.Range("E92").Value = 0.847383210396909
$ws.Range("E93").Value = 0.850895679662803
$ws.Range("E94").Value = 0.86406743940990505
$ws.Range("E95").Value = 0.82051282051282004
$ws.Range("E96").Value = 0.77994380049174505
$ws.Range("E97").Value = 0.87390235335440802

# ---------------------------------------------------------------------------
# 4. Highlight colours on column A of the summary table: the newly computed
#    rows are colour-coded (orange/green tones for top performers, red kept
#    for the rest).
# ---------------------------------------------------------------------------
$ws.Range("A90").Interior.Color = 7592334    # ~theme Accent6 tint 0.4 (teal/green)
$ws.Range("A97").Interior.Color = 7592334    # same colour, was red before

$ws.Range("A92").Interior.Color = 49407      # solid orange FFC000
$ws.Range("A93").Interior.Color = 49407
$ws.Range("A94").Interior.Color = 49407

$ws.Range("A95").Interior.Color = 255        # solid red FF0000 (newly highlighted)

# ---------------------------------------------------------------------------
# 5. Sheet view bookkeeping: re-centre the frozen-pane view & selection.
# ---------------------------------------------------------------------------
$excel.Goto($ws.Range("A77"), $true)
$ws.Range("C86").Select()
